# Incorporación de bees al Lector, Valuador y Sesionador
# Fill in the newly-valuated prices (column F) for the rows that previously
# showed "Sin precio".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$precios = @{
    2  = "907,38"
    4  = "907,38"
    5  = "942,72"
    6  = "1.248,26"
    7  = "1.384,14"
    8  = "1.310,41"
    11 = "676,99"
    12 = "676,99"
    23 = "232,90"
    28 = "1.671,37"
    30 = "1.671,37"
    37 = "552,99"
    39 = "815,64"
    40 = "815,64"
    42 = "393,16"
    43 = "552,99"
    46 = "393,16"
}

foreach ($fila in $precios.Keys) {
    $ws.Cells.Item($fila, 6).Value = $precios[$fila]
}
